# Updating rnasep2 IFB-Core assembly & annotation pipeline.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the assembly column headers ("Assembly_Raw"/"Assembly_Thin" -> "Raw"/"Thin").
$ws.Range("B1").Value = "Raw"
$ws.Range("C1").Value = "Thin"

# Fix the "Number of Trinity genes" raw-assembly figure (missing trailing digit).
$ws.Range("B2").Value = 331430

# Rename the "Basic" metric-type label to "Baseline" for the first block of rows.
$ws.Range("D2:D13").Value = "Baseline"

# Move the active selection to reflect where the author was last working.
$ws.Range("F10").Select()
